# Update "想去人数" (want-to-go count) values in column F
# for worksheets "展览" and "全部类型", matching the regenerated
# gh-pages data output at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Row -> new value for column F
$updates = @{
    2  = 1172
    3  = 596
    6  = 167
    8  = 63
    10 = 5441
    11 = 4863
    12 = 18
    15 = 52
    16 = 197
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
